$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 entirely (shift remaining rows up), leaving a 2-data-row table (rows 1-3)
$ws.Range("A4:T4").EntireRow.Delete() | Out-Null

# --- Update row 2 (previously FAPs/Rspo2/Lgr5/ECs) to FAPs/Rspo2/Lgr5/FAPs with recalculated TPM values ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3961209999999999
$ws.Range("N2").Value = 1.188363
$ws.Range("O2").Value = 0.5646784620538419
$ws.Range("P2").Value = 0.5646784620538419
$ws.Range("Q2").Value = 0.1377968957456666
$ws.Range("R2").Value = 1.240172061711
$ws.Range("S2").Value = 0.5646784620538419
$ws.Range("T2").Value = 0.5646784620538419

# --- Update row 3 (previously FAPs/Rspo2/Lgr5/FAPs) to FAPs/Rspo2/Lgr5/MuSCs with recalculated TPM values ---
$ws.Range("D3").Value = "MuSCs"
$ws.Range("M3").Value = 0.3053773333333333
$ws.Range("N3").Value = 0.9161319999999999
$ws.Range("O3").Value = 0.4353215379461581
$ws.Range("P3").Value = 0.4353215379461581
$ws.Range("Q3").Value = 0.1062302896448889
$ws.Range("R3").Value = 0.9560726068040001
$ws.Range("S3").Value = 0.4353215379461581
$ws.Range("T3").Value = 0.4353215379461581
